# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns with
# newly scraped values, and fix the ordering of the Cosmos / Injective
# Protocol rows (48/49) which swapped rank along with their data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $text) {
    # Force the cell to store a literal text value even when it looks like
    # a number (e.g. "512.83"), matching the source data where Price/Volume
    # columns are plain strings, not numerics.
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row => @{ D = new price (or $null to leave unchanged); E = new volume text }
$updates = @{
    2  = @{ D = "57.007.98";  E = "  +3.26%  " }
    3  = @{ D = "3.054.66";   E = "  +6.05%  " }
    4  = @{ D = $null;        E = "  -0.12%  " }
    5  = @{ D = "512.83";     E = "  +5.25%  " }
    6  = @{ D = "139.92";     E = "  +6.79%  " }
    7  = @{ D = "0.999";      E = "  -0.04%  " }
    8  = @{ D = $null;        E = "  +4.23%  " }
    9  = @{ D = "7.18";       E = "  +1.42%  " }
    10 = @{ D = "0.108";      E = "  +5.02%  " }
    11 = @{ D = "0.367";      E = "  +6.55%  " }
    12 = @{ D = "3.575.34";   E = "  +5.98%  " }
    13 = @{ D = $null;        E = "  +3.24%  " }
    14 = @{ D = "25.13";      E = "  -0.52%  " }
    15 = @{ D = $null;        E = "  +4.84%  " }
    16 = @{ D = "57.045.12";  E = "  +3.33%  " }
    17 = @{ D = "3.051.88";   E = "  +5.89%  " }
    18 = @{ D = $null;        E = "  -0.40%  " }
    19 = @{ D = "13.03";      E = "  +6.04%  " }
    20 = @{ D = "8.11";       E = "  +7.54%  " }
    21 = @{ D = "334.54";     E = "  +7.96%  " }
    22 = @{ D = "1.00";       E = "  +0.50%  " }
    23 = @{ D = "0.502";      E = "  +5.58%  " }
    24 = @{ D = "65.09";      E = "  +5.28%  " }
    25 = @{ D = $null;        E = "  +5.73%  " }
    26 = @{ D = $null;        E = "  +0.13%  " }
    27 = @{ D = "0.0₃0935";   E = "  +12.46%  " }
    28 = @{ D = "6.37";       E = "  +1.59%  " }
    29 = @{ D = "6.93";       E = "  +0.51%  " }
    30 = @{ D = "1.80";       E = "  +4.61%  " }
    31 = @{ D = "20.70";      E = "  +6.26%  " }
    32 = @{ D = $null;        E = "  +6.18%  " }
    33 = @{ D = "154.28";     E = "  +3.96%  " }
    34 = @{ D = "4.50";       E = "  +4.04%  " }
    35 = @{ D = "5.83";       E = "  +6.13%  " }
    36 = @{ D = "26.34";      E = "  +8.93%  " }
    37 = @{ D = $null;        E = "  +4.86%  " }
    38 = @{ D = "0.0667";     E = "  +3.61%  " }
    39 = @{ D = "3.088.12";   E = "  +6.02%  " }
    40 = @{ D = "36.84";      E = "  +2.45%  " }
    41 = @{ D = $null;        E = "  -0.14%  " }
    42 = @{ D = $null;        E = "  +6.62%  " }
    43 = @{ D = "3.80";       E = "  +5.68%  " }
    44 = @{ D = "2.226.07";   E = "  +6.97%  " }
    45 = @{ D = "0.0252";     E = "  +10.69%  " }
    46 = @{ D = $null;        E = "  +3.93%  " }
    47 = @{ D = "0.934";      E = "  +4.17%  " }
    50 = @{ D = "0.0860";     E = "  +3.18%  " }
    51 = @{ D = "0.683";      E = "  +6.56%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($null -ne $vals.D) {
        Set-TextValue $row 4 $vals.D
    }
    $ws.Cells.Item($row, 5).Value = $vals.E
}

# Rows 48 and 49 swapped places: Cosmos was rank 46 (row 48), Injective
# Protocol was rank 47 (row 49). Now Injective Protocol is rank 46 (row 48)
# and Cosmos is rank 47 (row 49), each with refreshed price/volume data.
$ws.Cells.Item(48, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue 48 4 "19.74"
$ws.Cells.Item(48, 5).Value = "  +7.17%  "

$ws.Cells.Item(49, 2).Value = "Cosmos"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue 49 4 "5.82"
$ws.Cells.Item(49, 5).Value = "  +0.70%  "
